{"js": "// Replace the three-digit-divided-by-one-digit equations in the table\n// with the new equations, preserving run formatting (font, size, etc.)\n// by using a search + in-place \"Replace\" insertText on each match.\n\nconst replacements = [\n  [\"238\u00f76=\", \"179\u00f72=\"],\n  [\"318\u00f79=\", \"845\u00f79=\"],\n  [\"816\u00f76=\", \"711\u00f78=\"],\n  [\"920\u00f72=\", \"380\u00f73=\"],\n  [\"858\u00f78=\", \"544\u00f76=\"],\n  [\"997\u00f75=\", \"113\u00f74=\"],\n  [\"422\u00f78=\", \"920\u00f79=\"],\n  [\"154\u00f79=\", \"356\u00f72=\"],\n  [\"613\u00f75=\", \"428\u00f75=\"],\n  [\"803\u00f79=\", \"943\u00f73=\"],\n  [\"734\u00f74=\", \"823\u00f78=\"],\n  [\"174\u00f76=\", \"156\u00f78=\"],\n  [\"811\u00f77=\", \"175\u00f79=\"],\n  [\"452\u00f75=\", \"991\u00f72=\"],\n  [\"567\u00f76=\", \"999\u00f78=\"],\n  [\"711\u00f74=\", \"504\u00f74=\"],\n  [\"195\u00f74=\", \"228\u00f73=\"],\n  [\"146\u00f79=\", \"594\u00f76=\"],\n  [\"802\u00f72=\", \"561\u00f78=\"],\n  [\"459\u00f79=\", \"590\u00f79=\"],\n  [\"963\u00f79=\", \"687\u00f77=\"],\n  [\"780\u00f76=\", \"184\u00f72=\"],\n  [\"578\u00f74=\", \"649\u00f74=\"],\n  [\"424\u00f79=\", \"102\u00f72=\"],\n  [\"302\u00f72=\", \"584\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"238\u00f76=\", \"179\u00f72=\"),\n    @(\"318\u00f79=\", \"845\u00f79=\"),\n    @(\"816\u00f76=\", \"711\u00f78=\"),\n    @(\"920\u00f72=\", \"380\u00f73=\"),\n    @(\"858\u00f78=\", \"544\u00f76=\"),\n    @(\"997\u00f75=\", \"113\u00f74=\"),\n    @(\"422\u00f78=\", \"920\u00f79=\"),\n    @(\"154\u00f79=\", \"356\u00f72=\"),\n    @(\"613\u00f75=\", \"428\u00f75=\"),\n    @(\"803\u00f79=\", \"943\u00f73=\"),\n    @(\"734\u00f74=\", \"823\u00f78=\"),\n    @(\"174\u00f76=\", \"156\u00f78=\"),\n    @(\"811\u00f77=\", \"175\u00f79=\"),\n    @(\"452\u00f75=\", \"991\u00f72=\"),\n    @(\"567\u00f76=\", \"999\u00f78=\"),\n    @(\"711\u00f74=\", \"504\u00f74=\"),\n    @(\"195\u00f74=\", \"228\u00f73=\"),\n    @(\"146\u00f79=\", \"594\u00f76=\"),\n    @(\"802\u00f72=\", \"561\u00f78=\"),\n    @(\"459\u00f79=\", \"590\u00f79=\"),\n    @(\"963\u00f79=\", \"687\u00f77=\"),\n    @(\"780\u00f76=\", \"184\u00f72=\"),\n    @(\"578\u00f74=\", \"649\u00f74=\"),\n    @(\"424\u00f79=\", \"102\u00f72=\"),\n    @(\"302\u00f72=\", \"584\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
